$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume update (GitHub Actions scheduled refresh).
# Some "Price" values are plain numeric-looking strings (e.g. "585.54") that
# Excel would otherwise auto-convert to numbers on assignment; force those
# specific cells to Text format first so they stay as strings, matching the
# original inline-string cell type used throughout column D.

$ws.Range('D2').Value = '64.287.71'
$ws.Range('E2').Value = '  +0.75%  '
$ws.Range('D3').Value = '3.488.42'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.54'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.27'
$ws.Range('E6').Value = '  +2.29%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -0.73%  '
$ws.Range('E9').Value = '  +0.71%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.21'
$ws.Range('E10').Value = '  +0.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.379'
$ws.Range('E11').Value = '  -1.09%  '
$ws.Range('D12').Value = '4.082.67'
$ws.Range('E12').Value = '  +0.09%  '
$ws.Range('E14').Value = '  +1.24%  '
$ws.Range('D15').Value = '3.488.68'
$ws.Range('E15').Value = '  +1.19%  '
$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '25.86'
$ws.Range('E16').Value = '  -6.18%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '64.297.65'
$ws.Range('E17').Value = '  +0.57%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.91'
$ws.Range('E18').Value = '  -2.42%  '
$ws.Range('E19').Value = '  +1.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.64'
$ws.Range('E20').Value = '  -5.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '387.63'
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('E22').Value = '  -2.10%  '
$ws.Range('D23').Value = '3.628.31'
$ws.Range('E23').Value = '  -0.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '74.27'
$ws.Range('E24').Value = '  +1.70%  '
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '5.72'
$ws.Range('E26').Value = '  +0.87%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.26'
$ws.Range('E30').Value = '  +0.63%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.49'
$ws.Range('E31').Value = '  -6.30%  '
$ws.Range('E32').Value = '  -0.61%  '
$ws.Range('D33').Value = '3.508.95'
$ws.Range('E33').Value = '  +0.33%  '
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E35').Value = '  +3.98%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.47'
$ws.Range('E36').Value = '  -1.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.20'
$ws.Range('E37').Value = '  -0.84%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.86'
$ws.Range('E38').Value = '  -0.05%  '
$ws.Range('E39').Value = '  -1.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '162.90'
$ws.Range('E40').Value = '  -2.71%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0782'
$ws.Range('E41').Value = '  -2.38%  '
$ws.Range('E42').Value = '  -0.82%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '25.40'
$ws.Range('E44').Value = '  -6.15%  '
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('E46').Value = '  +0.52%  '
$ws.Range('E47').Value = '  +1.46%  '
$ws.Range('E48').Value = '  -3.20%  '
$ws.Range('D49').Value = '2.473.96'
$ws.Range('E49').Value = '  +1.51%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.75'
$ws.Range('E50').Value = '  -1.71%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.902'
$ws.Range('E51').Value = '  +1.33%  '
